# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N ("Late"), shifting the existing "Late" / "Outstanding" / "heading"
# columns one place to the right (N->O, O->P, P->Q). The newly inserted
# column inherits the width of column M (the column immediately to its
# left), matching Excel's default "insert column" behaviour.
#
# Also make the "Repayment schedule" sheet the active sheet/tab, with
# K17 as the selected cell (previously "NewLoanInput" was the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of column M before inserting, so the new column can
# copy it (Excel's "insert" copies formatting/width from the column to
# the left when you right-click > Insert on a whole-column selection).
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column at N; everything from N onward shifts right.
$ws.Columns("N").Insert()

# New column N takes column M's width.
$ws.Columns("N").ColumnWidth = $leftWidth

# Switch to the Repayment schedule sheet and select K17, making this
# sheet/cell the active tab & selection for the workbook.
$ws.Range("K17").Select()
